$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the next day's date entry (continuing the work time table)
$null = $ws.Range("A7").Copy()
$null = $ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 43986

# Move active selection to B8, ready for the next hours entry
$null = $ws.Range("B8").Select()
